$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 5 values to the new set of test data strings
$ws.Range("A5").Value = "TestAutomation_8Sep"
$ws.Range("B5").Value = "A2241199332"
$ws.Range("C5").Value = "Facility_h224933552"
$ws.Range("D5").Value = "h224933552"
$ws.Range("E5").Value = "Pharmacy_h8Sept"
$ws.Range("F5").Value = "p199342232"
$ws.Range("G5").Value = "Cerner"
$ws.Range("H5").Value = "Alignment Project 192"

# Update the active selection on the sheet
$ws.Range("G9").Select()
